$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13; existing rows 13:22 shift down to 14:23,
# and the new blank row 13 inherits formatting from the row above (row 12 style).
$ws.Rows("13:13").Insert()

# Populate the newly inserted row 13 with the new record.
$ws.Range("A13").Value = 10
$ws.Range("B13").Value = "Vega Modelo de Temuco"
$ws.Range("C13").Value = "La Araucanía"
$ws.Range("D13").Value = 44777
$ws.Range("D13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E13").Value = 9
$ws.Range("F13").Value = 100112017
$ws.Range("G13").Value = "Ramas de apio"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 25
$ws.Range("K13").Value = 5000
$ws.Range("L13").Value = 5000
$ws.Range("M13").Value = 5000
$ws.Range("N13").Value = "$/paquete"
$ws.Range("O13").Value = "Región de La Araucanía"
$ws.Range("P13").Value = 5000
$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = "Hortaliza"
